# Reposition/resize the existing chart ("Chart 2") on the active sheet.
#
# Original anchor (twoCellAnchor): from col=10,colOff=0,row=7,rowOff=0
#                                   to   col=18,colOff=152400,row=22,rowOff=157163
# Target anchor:                    from col=4,colOff=504825,row=3,rowOff=152400
#                                   to   col=13,colOff=47625,row=19,rowOff=119063
#
# The sheet uses the default column width (58.4375 pt) and default row
# height (15 pt, per sheetFormatPr/defaultRowHeight), so the target
# top/left/width/height (in points) for the ChartObject are computed as:
#   left   = 4*58.4375 + 504825/12700  = 273.5
#   top    = 3*15      + 152400/12700  = 57
#   right  = 13*58.4375 + 47625/12700  = 763.4375
#   bottom = 19*15      + 119063/12700 = 294.37503937007875
#   width  = right - left  = 489.9375
#   height = bottom - top  = 237.37503937007875

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects("Chart 2")

$co.Left = 273.5
$co.Top = 57
$co.Width = 489.9375
$co.Height = 237.37503937007875
